# Daily status update - add entries for 13/11/2021 through 16/11/2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-StatusCell {
    param($addr, $text, $dateCol)

    $rng = $ws.Range($addr)
    $rng.Value2 = $text
    if ($dateCol) {
        $rng.NumberFormat = "@"
    }
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
}

# Row 105 - 13/11/2021 (Holiday)
Set-StatusCell "A105" "13/11/2021" $true
Set-StatusCell "B105" "HOLIDAY" $false

# Row 106 - 14/11/2021 (Holiday)
Set-StatusCell "A106" "14/11/2021" $true
Set-StatusCell "B106" "HOLIDAY" $false

# Row 107 - 15/11/2021
Set-StatusCell "A107" "15/11/2021" $true
Set-StatusCell "B107" "OpenMax : Resource management,Buffer Payload" $false

# Row 108 - additional notes for 15/11/2021
Set-StatusCell "B108" "Internal discussion with teammates : Testapps" $false
Set-StatusCell "C108" "Updating the notes" $false
Set-StatusCell "D108" "Revision of previous topics" $false

# Row 109 - additional note
Set-StatusCell "B109" "LDD recap session" $false

# Row 110 - 16/11/2021
Set-StatusCell "A110" "16/11/2021" $true

# Update the view to match where the user left off editing
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 99
$aw.ScrollColumn = 1
$ws.Range("A110").Select()
